$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values entered in the same order the author typed them, so the
# shared-string table comes out in the same sequence as the target file.
$ws.Range("C4").Value = "melvinngweixiang@gmail.com,2101000d@student.tp.edu.sg"
$ws.Range("D4").Value = "Test3"
$ws.Range("D5").Value = "Test4"
$ws.Range("E5").Value = "Testing4"
$ws.Range("E4").Value = "Testing3"

$ws.Range("A4").Value = "Gmail"
$ws.Range("B4").Value = "melvinngweixiang@gmail.com"
$ws.Range("F4").Value = "No"

$ws.Range("A5").Value = "Outlook"
$ws.Range("B5").Value = "melvinngweixiang@gmail.com"
$ws.Range("C5").Value = "melvinngweixiang@gmail.com,2101000d@student.tp.edu.sg"
$ws.Range("F5").Value = "No"

# Hyperlinks, in the same order they appear in the target file (B4, B5, C4, C5)
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:melvinngweixiang@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:melvinngweixiang@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:melvinngweixiang@gmail.com,2101000d@student.tp.edu.sg") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:melvinngweixiang@gmail.com,2101000d@student.tp.edu.sg") | Out-Null

# Hyperlinks.Add stamps a fresh "applied hyperlink" style; re-apply the
# named Hyperlink style afterwards so these cells share style index 1
# with the pre-existing hyperlink cells (B2/B3/C3) instead of a new one.
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("C4").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("C5").Style = "Hyperlink"

# Row 6 - only B6 carries hyperlink style, no value
$ws.Range("B6").Style = "Hyperlink"

$ws.Range("C10").Select()
